$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 9.2.1 indicator text reworded
$ws.Range("B4").Value = "9.2.1 Добавленная стоимость, создаваемая в обрабатывающей промышленности, в процентном отношении к ВВП и на душу населения"

# Organization website updated
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Leave the last active cell on B4, matching the saved selection in the workbook
$ws.Range("B4").Select()
